$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.693.89"
$ws.Range("E2").Value = "  +2.29%  "

$ws.Range("D3").Value = "2.200.97"
$ws.Range("E3").Value = "  -0.49%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "257.58"
$ws.Range("E5").Value = "  +1.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "84.51"
$ws.Range("E6").Value = "  +11.98%  "

$ws.Range("E7").Value = "  +0.75%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  +1.56%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.82"
$ws.Range("E10").Value = "  +9.57%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0922"
$ws.Range("E11").Value = "  +0.11%  "

$ws.Range("E12").Value = "  +4.62%  "

$ws.Range("E13").Value = "  +2.43%  "

$ws.Range("D14").Value = "2.531.37"
$ws.Range("E14").Value = "  -0.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.38"
$ws.Range("E15").Value = "  +0.39%  "

$ws.Range("D16").Value = "2.206.60"
$ws.Range("E16").Value = "  -0.05%  "

$ws.Range("E17").Value = "  +0.35%  "

$ws.Range("D18").Value = "43.638.82"
$ws.Range("E18").Value = "  +2.36%  "

$ws.Range("E19").Value = "  +0.11%  "

$ws.Range("E20").Value = "  -1.72%  "

$ws.Range("E21").Value = "  -0.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.37"
$ws.Range("E22").Value = "  +7.97%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.57"
$ws.Range("E23").Value = "  +1.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.14"
$ws.Range("E24").Value = "  -3.07%  "

$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.63"
$ws.Range("E26").Value = "  +8.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.62"
$ws.Range("E27").Value = "  +0.44%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.21"
$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.26"
$ws.Range("E29").Value = "  +3.10%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.23"
$ws.Range("E30").Value = "  +2.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.65"
$ws.Range("E31").Value = "  +0.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.41"
$ws.Range("E32").Value = "  +0.92%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0861"
$ws.Range("E33").Value = "  +3.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.36"
$ws.Range("E34").Value = "  +3.12%  "

$ws.Range("E35").Value = "  +1.78%  "

$ws.Range("E36").Value = "  +2.05%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.50"
$ws.Range("E37").Value = "  +5.23%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0360"
$ws.Range("E38").Value = "  +4.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.46"
$ws.Range("E39").Value = "  +1.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.84"
$ws.Range("E40").Value = "  +5.23%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.10"
$ws.Range("E41").Value = "  +0.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "62.95"
$ws.Range("E42").Value = "  +5.09%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.48"
$ws.Range("E43").Value = "  +4.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.199"
$ws.Range("E44").Value = "  +1.75%  "

$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0982"
$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.17"
$ws.Range("E46").Value = "  -2.61%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.32"
$ws.Range("E47").Value = "  -0.56%  "

$ws.Range("E48").Value = "  +5.20%  "

$ws.Range("E49").Value = "  +0.55%  "

$ws.Range("E50").Value = "  -5.36%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.50"
$ws.Range("E51").Value = "  +8.56%  "
